# Applies the data refresh described in the commit:
#   "Update gh-pages to output generated at 456a3b4"
# Numeric "want-to-go" counters bump slightly and a handful of
# rows on the "全部类型" (all types) sheet are replaced with newer events.
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(8,6).Value = 30
$ws.Cells.Item(11,6).Value = 1281
$ws.Cells.Item(12,6).Value = 29440
$ws.Cells.Item(13,6).Value = 4897
$ws.Cells.Item(17,6).Value = 57
$ws.Cells.Item(19,6).Value = 26
$ws.Cells.Item(22,6).Value = 18
$ws.Cells.Item(23,6).Value = 653
$ws.Cells.Item(24,6).Value = 283
$ws.Cells.Item(28,6).Value = 89
$ws.Cells.Item(29,6).Value = 10
$ws.Cells.Item(30,6).Value = 677
$ws.Cells.Item(33,6).Value = 566
$ws.Cells.Item(36,6).Value = 656
# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(6,6).Value = 391
$ws.Cells.Item(7,6).Value = 923
$ws.Cells.Item(12,6).Value = 4270
$ws.Cells.Item(23,6).Value = 4256
# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3,6).Value = 272
$ws.Cells.Item(4,6).Value = 1235
$ws.Cells.Item(5,6).Value = 309
# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3,6).Value = 272
$ws.Cells.Item(4,6).Value = 1235
$ws.Cells.Item(7,6).Value = 391
$ws.Cells.Item(8,6).Value = 309
$ws.Cells.Item(9,6).Value = 923
# Row 11
$ws.Cells.Item(11,3).Value = '广州·《龙珠》《灌篮高手》回忆与幻想——世界经典动漫主题音乐会'
$ws.Cells.Item(11,4).Value = '天河路228号正佳广场七楼 广州正佳大剧院'
$ws.Cells.Item(11,5).Value = '2024.07.06 15:30-07.06 17:00'
$ws.Cells.Item(11,6).Value = 2
$ws.Cells.Item(11,7).Value = 50
$ws.Cells.Item(11,8).Value = 'https://show.bilibili.com/platform/detail.html?id=88096'
$ws.Cells.Item(11,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/xYauX6km1719314524354.jpeg'
$ws.Cells.Item(14,6).Value = 30
$ws.Cells.Item(18,6).Value = 1281
# Row 22
$ws.Cells.Item(22,3).Value = '广州·燃动!!高梨康治SUMMER LIVE TOUR IN CHINA 2024'
$ws.Cells.Item(22,4).Value = '海珠同创汇东一街11号（上冲南约11-2） 声音共和Livehouse'
$ws.Cells.Item(22,5).Value = '2024.07.21 14:30-07.21 16:00'
$ws.Cells.Item(22,6).Value = 192
$ws.Cells.Item(22,7).Value = 280
$ws.Cells.Item(22,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87034'
$ws.Cells.Item(22,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/LINsP2ui1717741701901.png'
# Row 24
$ws.Cells.Item(24,2).NumberFormat = "@"
$ws.Cells.Item(24,2).Value = '2024-07-26'
$ws.Cells.Item(24,2).Style = "Normal"
$ws.Cells.Item(24,3).Value = '广州·【早鸟8折】“浪漫古典Ⅱ”百年经典传世名曲烛光音乐会 '
$ws.Cells.Item(24,4).Value = '广州市二沙岛晴波路33号  星海音乐厅（交响乐演奏厅）'
$ws.Cells.Item(24,5).Value = '2024.07.26 20:00-07.26 21:30'
$ws.Cells.Item(24,6).Value = 1
$ws.Cells.Item(24,7).Value = 144
$ws.Cells.Item(24,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87726'
$ws.Cells.Item(24,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/A8vhVlhn1717575084179.png'
# Row 25
$ws.Cells.Item(25,3).Value = '广州·萨克斯王子安德鲁·杨——2024经典&流行音乐巡回演出'
$ws.Cells.Item(25,4).Value = '龙凤街道革新路124号太古仓码头5号仓 广州太空间Live House'
$ws.Cells.Item(25,6).Value = 4
$ws.Cells.Item(25,7).Value = 380
$ws.Cells.Item(25,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86635'
$ws.Cells.Item(25,9).Value = '//i1.hdslb.com/bfs/openplatform/202405/rciNih361716802006584.jpeg'
# Row 26
$ws.Cells.Item(26,2).NumberFormat = "@"
$ws.Cells.Item(26,2).Value = '2024-07-27'
$ws.Cells.Item(26,2).Style = "Normal"
$ws.Cells.Item(26,3).Value = '广州·LookLook动漫嘉年华2th'
$ws.Cells.Item(26,4).Value = '东沙大道16号 健康方舟6层博览馆'
$ws.Cells.Item(26,5).Value = '2024.07.27 10:00-07.28 17:30'
$ws.Cells.Item(26,6).Value = 270
$ws.Cells.Item(26,7).Value = 68
$ws.Cells.Item(26,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87217'
$ws.Cells.Item(26,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/p4Bs2soo1718113055492.png'
$ws.Cells.Item(30,6).Value = 57
$ws.Cells.Item(31,6).Value = 26
# Row 35
$ws.Cells.Item(35,3).Value = '广州·代号鸢only'
$ws.Cells.Item(35,4).Value = '会江路 巨大产业园智慧港'
$ws.Cells.Item(35,5).Value = '2024.08.03 10:00-08.04 17:00'
$ws.Cells.Item(35,6).Value = 18
$ws.Cells.Item(35,7).Value = 55
$ws.Cells.Item(35,8).Value = 'https://show.bilibili.com/platform/detail.html?id=88224'
$ws.Cells.Item(35,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/qBntv1WK1719481529863.jpeg'
$ws.Cells.Item(36,6).Value = 653
$ws.Cells.Item(37,6).Value = 283
$ws.Cells.Item(39,6).Value = 89
$ws.Cells.Item(40,6).Value = 10
$ws.Cells.Item(41,6).Value = 677
$ws.Cells.Item(49,6).Value = 656
